$d = $word.ActiveDocument

# 1. Skills line: add ", SFCC/SFRA" to the JS/jQuery/Node/React line.
$d.Content.Find.Execute("JS, jQuery, Node, React, ES6 modular development, PHP", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "JS, jQuery, Node, React, ES6 modular development, PHP, SFCC/SFRA", 2)

# 2. Insert three new bullet points right before the "Redesigned and updated PDP elements..." bullet
#    (currently the 26th paragraph in the document).
$pdpPara = $d.Paragraphs(26)
$pdpPara.Range.InsertBefore("Rebuilt the Quantity Stepper from a dropodown to an input with plus/minus buttons that supported a number of business requirements including but not limited to: keyboard and accessibility, available stock including the number left if under threshold defined in Business Manager, restrict number of items to purchase including those already in cart (also defined in Business Manager), and would display an error message anytime the parameters were exceeded.`rUtilizing Business Manager, Customs Jobs and Objects built the Back In Stock Notifications that would add the user and product to a custom object for a custom job to loop through once a day and check if products were back in stock and if so, would add the user and details to our email client to be notified.`rUtilizing Business Manager and content assets built a Size Comparison Chart modal that would pull in any content asset to display to the user.`r")

# 3. Reviewed project requirements ... add "(AGILE)"
$d.Content.Find.Execute("Reviewed project requirements to align with merchandiser to stay within resource constraints.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Reviewed project requirements to align with merchandiser to stay within resource constraints (AGILE).", 2)

# 4. "can manage" -> "too manage"
$d.Content.Find.Execute([string]::Concat("Merchandiser", [char]0x2019, "s/Content Manager", [char]0x2019, "s can manage each individual storefront/language."), $true, $false, $false, $false, $false, `
                         $true, 1, $false, [string]::Concat("Merchandiser", [char]0x2019, "s/Content Manager", [char]0x2019, "s too manage each individual storefront/language."), 2)

# 5. "European audiences" -> "Worldwide audiences." (merges trailing "." run into the text run and
#    drops the "European"/"Worldwide" wording change), then relocate the _GoBack bookmark that used
#    to sit between the two runs into the following (blank) paragraph.
$d.Content.Find.Execute("Assisted in localizing entire marketing site for our European audiences.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Assisted in localizing entire marketing site for our Worldwide audiences.", 2)

$blankPara = $d.Paragraphs(47)
$blankRange = $blankPara.Range
$blankRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $blankRange)
